$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.100.86"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -3.54%  '

$ws.Range("D3").Value = "'1.651.85"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -5.31%  '

$ws.Range("D4").Value = "'0.9998"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").Value = "'237.80"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.50%  '

$ws.Range("E6").Value = '  +0.05%  '

$ws.Range("D7").Value = "'0.4798"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -7.08%  '

$ws.Range("D8").Value = "'0.2617"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -5.23%  '

$ws.Range("D9").Value = "'0.06011"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.06%  '

$ws.Range("D10").Value = "'0.07188"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").Value = "'1.659.51"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -4.89%  '

$ws.Range("D12").Value = "'14.86"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -2.55%  '

$ws.Range("D13").Value = "'0.6242"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.95%  '

$ws.Range("D14").Value = "'4.609"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -0.64%  '

$ws.Range("D15").Value = "'73.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -5.73%  '

$ws.Range("D16").Value = "'0.9999"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +0.00%  '

$ws.Range("D17").Value = "'0.9994"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -0.08%  '

$ws.Range("D18").Value = "'25.089.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.70%  '

$ws.Range("D19").Value = "'11.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -3.56%  '

$ws.Range("D20").Value = "'0.000006610"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -3.01%  '

$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").Value = "'4.473"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +4.25%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = "'1.862.52"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -5.32%  '

$ws.Range("D23").Value = "'8.628"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -0.66%  '

$ws.Range("D24").Value = "'5.296"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -1.13%  '

$ws.Range("D25").Value = "'133.18"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.87%  '

$ws.Range("D26").Value = "'14.97"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.17%  '

$ws.Range("D27").Value = "'1.393"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -7.63%  '

$ws.Range("D28").Value = "'103.58"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -2.04%  '

$ws.Range("D29").Value = "'1.684"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -5.51%  '

$ws.Range("D30").Value = "'3.781"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -4.53%  '

$ws.Range("D31").Value = "'0.07919"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -4.42%  '

$ws.Range("D32").Value = "'3.576"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.78%  '

$ws.Range("D33").Value = "'0.04593"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.93%  '

$ws.Range("D34").Value = "'2.582"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -2.81%  '

$ws.Range("D35").Value = "'0.9459"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -5.67%  '

$ws.Range("D36").Value = "'0.5774"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -7.61%  '

$ws.Range("D37").Value = "'2.621"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -4.19%  '

$ws.Range("D38").Value = "'0.01555"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -3.07%  '

$ws.Range("D39").Value = "'0.8396"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  +11.07%  '

$ws.Range("E40").Value = '  +0.09%  '

$ws.Range("D41").Value = "'1.832"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -5.42%  '

$ws.Range("D42").Value = "'99.10"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.47%  '

$ws.Range("D43").Value = "'0.3721"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -4.12%  '

$ws.Range("D44").Value = "'4.809"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -4.14%  '

$ws.Range("D45").Value = "'0.1137"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.13%  '

$ws.Range("D46").Value = "'6.125"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -3.53%  '

$ws.Range("E47").Value = '  -0.91%  '

$ws.Range("D48").Value = "'29.85"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -2.68%  '

$ws.Range("D49").Value = "'51.47"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -7.35%  '

$ws.Range("E50").Value = '  -0.01%  '

$ws.Range("B51").Value = 'Decentraland'
$ws.Range("C51").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D51").Value = "'0.3343"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -2.95%  '
